$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion note text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
[string]$text = $cellA1.Value()
$text = $text.Replace("1000 Bs = 7.07 = 28091.87 pesos", "1000 Bs = 7.05 = 28141.13 pesos")
$text = $text.Replace("28091.87 pesos = 7.02 = 963.58 Bs", "28141.13 pesos = 7.02 = 971.04 Bs")
$cellA1.Value = $text

# --- Sheet "tasas": update the rate figures in N10/O10 and N12/O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 141.86
$wsTasas.Range("O10").Value = 3992.1
$wsTasas.Range("N12").Value = 4008
$wsTasas.Range("O12").Value = 138.3
